# Applies the annotation-sample additions described in the commit:
# "Add annotation files and update datasets" — specifically the new
# annot3_* sample values added to the hypercap_cohort sheet (columns
# BM:BP, rows 2-10), plus the corresponding selection updates on the
# hypercap_cohort and Data sheets.

$wb = $excel.ActiveWorkbook

$wsCohort = $wb.Worksheets.Item("hypercap_cohort")
$wsData   = $wb.Worksheets.Item("Data")

# ---------------------------------------------------------------------
# 1. Write the new annotation sample values into BM2:BP10.
#    (Text matches the existing "Data" sheet / header vocabulary so the
#    shared-string table is reused rather than duplicated.)
# ---------------------------------------------------------------------

$wsCohort.Range("BM2").Value = "Chest pain"
$wsCohort.Range("BN2").Value = "Symptom – Circulatory"
$wsCohort.Range("BO2").Value = "Dyspnea"
$wsCohort.Range("BP2").Value = "Symptom – Respiratory"

$wsCohort.Range("BM3").Value = "Altered mental status"
$wsCohort.Range("BN3").Value = "Symptom – Nervous"
$wsCohort.Range("BO3").Value = "Weakness"
$wsCohort.Range("BP3").Value = "Symptom – General"

$wsCohort.Range("BM4").Value = "Dyspnea, Respiratory distress"
$wsCohort.Range("BN4").Value = "Symptom – Respiratory"

$wsCohort.Range("BM5").Value = "Weakness"
$wsCohort.Range("BN5").Value = "Symptom – General"

$wsCohort.Range("BM6").Value = "HYPOTENSIVE"
$wsCohort.Range("BN6").Value = "Symptom – Circulatory"

$wsCohort.Range("BM7").Value = "Overdose"
$wsCohort.Range("BN7").Value = "Injuries & adverse effects"

$wsCohort.Range("BM8").Value = "Slurred speech"
$wsCohort.Range("BN8").Value = "Symptom – Nervous"

$wsCohort.Range("BM9").Value = "Transfer, Respiratory distress"

$wsCohort.Range("BM10").Value = "Dyspnea"
$wsCohort.Range("BN10").Value = "Symptom – Respiratory"

# ---------------------------------------------------------------------
# 2. Re-apply the shaded "highlight" look (same fill already used on the
#    chief_complaint column, S) to the BM column for rows 3-10, by
#    copying the format from an already-shaded cell.
# ---------------------------------------------------------------------

$wsCohort.Range("S2").Copy() | Out-Null
$wsCohort.Range("BM3:BM10").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 3. Give the BN5, BN7 and BP3 "cat" cells the new Cambria font
#    treatment. A temporary named style is used (and removed again) so
#    that exactly one new font and one new cell format are added to the
#    workbook, matching what Excel itself produces for this kind of
#    direct formatting.
# ---------------------------------------------------------------------

$cambriaStyle = $wb.Styles.Add("CambriaTemp")
$cambriaStyle.Font.Name = "Cambria"

$wsCohort.Range("BP3").Style = "CambriaTemp"
$wsCohort.Range("BN5").Style = "CambriaTemp"
$wsCohort.Range("BN7").Style = "CambriaTemp"

$cambriaStyle.Delete()

# ---------------------------------------------------------------------
# 4. Update the on-screen selection / scroll position to match the
#    author's final view of the hypercap_cohort sheet.
# ---------------------------------------------------------------------

$wsCohort.Activate()
$wsCohort.Range("BN11").Select()

# ---------------------------------------------------------------------
# 5. Update the selection on the Data sheet as well.
# ---------------------------------------------------------------------

$wsData.Activate()
$wsData.Range("A8").Select()

# Leave the cohort sheet as the active/visible tab, matching tabSelected.
$wsCohort.Activate()
